$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '36.790.46'
Set-TextValue 2 5 '  -1.14%  '
Set-TextValue 3 4 '2.092.08'
Set-TextValue 3 5 '  +1.73%  '
Set-TextValue 4 5 '  +0.09%  '
Set-TextValue 5 4 '245.40'
Set-TextValue 5 5 '  -1.35%  '
Set-TextValue 6 5 '  -1.75%  '
Set-TextValue 8 4 '54.48'
Set-TextValue 8 5 '  -6.10%  '
Set-TextValue 9 4 '58.85'
Set-TextValue 9 5 '  -2.15%  '
Set-TextValue 10 5 '  -4.27%  '
Set-TextValue 11 5 '  -2.28%  '
Set-TextValue 12 5 '  +0.84%  '
Set-TextValue 13 4 '0.907'
Set-TextValue 13 5 '  +3.24%  '
Set-TextValue 14 4 '15.09'
Set-TextValue 14 5 '  -5.83%  '
Set-TextValue 15 4 '2.395.61'
Set-TextValue 15 5 '  +1.76%  '
Set-TextValue 16 5 '  -3.58%  '
Set-TextValue 17 4 '2.102.60'
Set-TextValue 17 5 '  +2.32%  '
Set-TextValue 18 4 '36.764.86'
Set-TextValue 18 5 '  -1.09%  '
Set-TextValue 19 4 '17.15'
Set-TextValue 19 5 '  -6.60%  '
Set-TextValue 20 4 '72.73'
Set-TextValue 20 5 '  -2.97%  '
Set-TextValue 21 5 '  -1.26%  '
Set-TextValue 22 5 '  +0.84%  '
Set-TextValue 23 4 '238.81'
Set-TextValue 23 5 '  +0.54%  '
Set-TextValue 24 5 '  +0.08%  '
Set-TextValue 25 5 '  -3.20%  '
Set-TextValue 26 4 '9.74'
Set-TextValue 26 5 '  +2.32%  '
Set-TextValue 27 5 '  -1.16%  '
Set-TextValue 28 4 '167.49'
Set-TextValue 28 5 '  -1.25%  '
Set-TextValue 29 5 '  +2.36%  '
Set-TextValue 30 5 '  -1.56%  '
Set-TextValue 31 5 '  +9.50%  '
Set-TextValue 32 4 '1.17'
Set-TextValue 32 5 '  +3.15%  '
Set-TextValue 33 4 '4.74'
Set-TextValue 33 5 '  +5.48%  '
Set-TextValue 34 4 '0.0611'
Set-TextValue 34 5 '  -1.61%  '
Set-TextValue 35 4 '2.43'
Set-TextValue 35 5 '  +7.55%  '
Set-TextValue 36 5 '  +0.23%  '
Set-TextValue 37 5 '  +3.82%  '
Set-TextValue 38 4 '0.0833'
Set-TextValue 38 5 '  -6.73%  '
Set-TextValue 40 5 '  +1.07%  '
Set-TextValue 42 4 '4.87'
Set-TextValue 42 5 '  -7.70%  '
Set-TextValue 43 4 '0.0955'
Set-TextValue 43 5 '  -3.11%  '
Set-TextValue 44 4 '96.40'
Set-TextValue 44 5 '  +0.05%  '
Set-TextValue 45 5 '  -9.50%  '
Set-TextValue 46 4 '16.10'
Set-TextValue 46 5 '  -7.23%  '
Set-TextValue 47 4 '1.382.46'
Set-TextValue 47 5 '  +8.66%  '
Set-TextValue 48 4 '7.41'
Set-TextValue 48 5 '  +8.42%  '
Set-TextValue 49 5 '  +0.34%  '
Set-TextValue 50 4 '2.90'
Set-TextValue 50 5 '  +1.41%  '
Set-TextValue 51 4 '2.282.77'
Set-TextValue 51 5 '  +1.92%  '
